# Generate Report for Handoff
# Updates the localization-status report: the first source file ("9c3d37a6-...")
# was handed off again (new guid "bd32cb5f-..."), and the second source file
# ("f1023a51-...") was detected as a content duplicate of the first, so it now
# shares the same handoff artifact and is renamed to a new guid
# ("ffffb926d66b-...").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$oldGuid1 = "9c3d37a6-6f74-4991-b88d-a574185a4985"
$newGuid1 = "bd32cb5f-a5de-46e6-b32d-cc1b6cc08662"
$oldGuid2 = "f1023a51-e6f7-44c0-8916-f06240439146"
$newGuid2 = "ffffb926d66b-38dc-4864-8669-61dfd909342b"

$newStatus = "Ready for handoff"
$newHoDate = "2016-08-25 17:05:20"
$newHandoffDateZh = "2016-08-25 17:05:15"
$newHandoffDateDe = "2016-08-25 17:05:20"
$emptyHandback = "0001-01-01 00:00:00"
$newHash = "be67c59b958f36f60f42e30976160bd87187b347"

$newZhXlf = "$newGuid1.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid1.$newHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "$newGuid1.md"
$ws1.Range("B2").Value = "e2e\$newGuid1.md"
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("G2").Value = $newHoDate

$ws1.Range("A3").Value = "$newGuid2.md"
$ws1.Range("B3").Value = "e2e\$newGuid2.md"
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus
$ws1.Range("G3").Value = $newHoDate

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3599633a9b06ce9cdf803055ddaa90dd177d745/e2e/$oldGuid1.md", "", "", "e2e\$newGuid1.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3599633a9b06ce9cdf803055ddaa90dd177d745/e2e/$oldGuid2.md", "", "", "e2e\$newGuid2.md") | Out-Null

$ws1.Columns.Item(5).AutoFit() | Out-Null
$ws1.Columns.Item(6).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "$newGuid1.md"
$ws2.Range("C2").Value = $newStatus
$ws2.Range("G2").Value = $newZhXlf
$ws2.Range("H2").Value = $newHandoffDateZh
$ws2.Range("I2").Value = ""
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = $emptyHandback

$ws2.Range("A3").Value = "$newGuid2.md"
$ws2.Range("C3").Value = $newStatus
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = $newZhXlf
$ws2.Range("H3").Value = $newHandoffDateZh
$ws2.Range("I3").Value = ""
$ws2.Range("J3").Value = ""
$ws2.Range("K3").Value = $emptyHandback

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3599633a9b06ce9cdf803055ddaa90dd177d745/e2e/$oldGuid1.md", "", "", "$newGuid1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3599633a9b06ce9cdf803055ddaa90dd177d745/e2e/$oldGuid2.md", "", "", "$newGuid2.md") | Out-Null

$ws2.Columns.Item(3).AutoFit() | Out-Null
$ws2.Columns.Item(9).AutoFit() | Out-Null
$ws2.Columns.Item(10).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = "$newGuid1.md"
$ws3.Range("C2").Value = $newStatus
$ws3.Range("G2").Value = $newDeXlf
$ws3.Range("H2").Value = $newHandoffDateDe
$ws3.Range("I2").Value = ""
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = $emptyHandback

$ws3.Range("A3").Value = "$newGuid2.md"
$ws3.Range("C3").Value = $newStatus
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = $newDeXlf
$ws3.Range("H3").Value = $newHandoffDateDe
$ws3.Range("I3").Value = ""
$ws3.Range("J3").Value = ""
$ws3.Range("K3").Value = $emptyHandback

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3599633a9b06ce9cdf803055ddaa90dd177d745/e2e/$oldGuid1.md", "", "", "$newGuid1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3599633a9b06ce9cdf803055ddaa90dd177d745/e2e/$oldGuid2.md", "", "", "$newGuid2.md") | Out-Null

$ws3.Columns.Item(3).AutoFit() | Out-Null
$ws3.Columns.Item(9).AutoFit() | Out-Null
$ws3.Columns.Item(10).AutoFit() | Out-Null

Write-Host "Report updated for handoff."
